$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") from 45192 -> 45202 for the existing data rows (2..511)
for ($r = 2; $r -le 511; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}

# 2) Row 511 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(511).RowHeight = 15

# 3) Append three new data rows: 512, 513, 514

# --- Row 512 ---
$ws.Cells.Item(512, 1).Value = "A 45829-2023"
$ws.Cells.Item(512, 2).Value = 45195
$ws.Cells.Item(512, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(512, 3).Value = 45202
$ws.Cells.Item(512, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(512, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(512, 5).Value = "FALKENBERG"
$ws.Cells.Item(512, 7).Value = 1.8
$ws.Cells.Item(512, 8).Value = 0
$ws.Cells.Item(512, 9).Value = 0
$ws.Cells.Item(512, 10).Value = 0
$ws.Cells.Item(512, 11).Value = 0
$ws.Cells.Item(512, 12).Value = 0
$ws.Cells.Item(512, 13).Value = 0
$ws.Cells.Item(512, 14).Value = 0
$ws.Cells.Item(512, 15).Value = 0
$ws.Cells.Item(512, 16).Value = 0
$ws.Cells.Item(512, 17).Value = 0
$ws.Cells.Item(512, 18).WrapText = $true
$ws.Rows.Item(512).RowHeight = 15

# --- Row 513 ---
$ws.Cells.Item(513, 1).Value = "A 46059-2023"
$ws.Cells.Item(513, 2).Value = 45196
$ws.Cells.Item(513, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(513, 3).Value = 45202
$ws.Cells.Item(513, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(513, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(513, 5).Value = "FALKENBERG"
$ws.Cells.Item(513, 7).Value = 1.7
$ws.Cells.Item(513, 8).Value = 0
$ws.Cells.Item(513, 9).Value = 0
$ws.Cells.Item(513, 10).Value = 0
$ws.Cells.Item(513, 11).Value = 0
$ws.Cells.Item(513, 12).Value = 0
$ws.Cells.Item(513, 13).Value = 0
$ws.Cells.Item(513, 14).Value = 0
$ws.Cells.Item(513, 15).Value = 0
$ws.Cells.Item(513, 16).Value = 0
$ws.Cells.Item(513, 17).Value = 0
$ws.Cells.Item(513, 18).WrapText = $true
$ws.Rows.Item(513).RowHeight = 15

# --- Row 514 ---
$ws.Cells.Item(514, 1).Value = "A 46924-2023"
$ws.Cells.Item(514, 2).Value = 45201
$ws.Cells.Item(514, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(514, 3).Value = 45202
$ws.Cells.Item(514, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(514, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(514, 5).Value = "FALKENBERG"
$ws.Cells.Item(514, 6).Value = "Kyrkan"
$ws.Cells.Item(514, 7).Value = 4
$ws.Cells.Item(514, 8).Value = 0
$ws.Cells.Item(514, 9).Value = 0
$ws.Cells.Item(514, 10).Value = 0
$ws.Cells.Item(514, 11).Value = 0
$ws.Cells.Item(514, 12).Value = 0
$ws.Cells.Item(514, 13).Value = 0
$ws.Cells.Item(514, 14).Value = 0
$ws.Cells.Item(514, 15).Value = 0
$ws.Cells.Item(514, 16).Value = 0
$ws.Cells.Item(514, 17).Value = 0
$ws.Cells.Item(514, 18).WrapText = $true

Write-Output "done"
